$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3670.182
$ws.Range("J125").Value = 3670.182
$ws.Range("L125").Value = 33031.638
$ws.Range("N125").Value = -37951.638
$ws.Range("H129").Value = 870.7474999999999
$ws.Range("I129").Value = 463.89474
$ws.Range("J129").Value = 967.375
$ws.Range("K129").Value = 1391.68422
$ws.Range("L129").Value = 2902.125
$ws.Range("M129").Value = 3608.31578
$ws.Range("N129").Value = -12902.125
$ws.Range("H137").Value = 1445.25
$ws.Range("I137").Value = 1378.5454
$ws.Range("K137").Value = 4135.6362
$ws.Range("M137").Value = -1585.6362
$ws.Range("H138").Value = 2751.3289
$ws.Range("I138").Value = 726.08887
$ws.Range("J138").Value = 5691.1934
$ws.Range("K138").Value = 2178.26661
$ws.Range("L138").Value = 17073.5802
$ws.Range("M138").Value = 2961.73339
$ws.Range("N138").Value = -27353.5802
$ws.Range("H141").Value = 1149.55
$ws.Range("I141").Value = 1104.7894
$ws.Range("K141").Value = 3314.3682
$ws.Range("M141").Value = 1865.6318

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 15050
$ws.Range("J44").Value = 15050
$ws.Range("L44").Value = 15050
$ws.Range("N44").Value = -16044
$ws.Range("H134").Value = 6665.88
$ws.Range("I134").Value = 8764.5625
$ws.Range("J134").Value = 2934.889
$ws.Range("K134").Value = 26293.6875
$ws.Range("L134").Value = 8804.667000000001
$ws.Range("M134").Value = -23758.6875
$ws.Range("N134").Value = -13874.667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11908232
$ws.Range("I31").Value = 1546.44
$ws.Range("J31").Value = 29418062
$ws.Range("K31").Value = 1546.44
$ws.Range("L31").Value = 29418062
$ws.Range("M31").Value = -1251.44
$ws.Range("N31").Value = -29418652
$ws.Range("H34").Value = 11908232
$ws.Range("I34").Value = 1546.44
$ws.Range("J34").Value = 29418062
$ws.Range("K34").Value = 1546.44
$ws.Range("L34").Value = 29418062
$ws.Range("M34").Value = -1344.44
$ws.Range("N34").Value = -29418466
$ws.Range("H58").Value = 5209702
$ws.Range("I58").Value = 7576701
$ws.Range("J58").Value = 2303.6
$ws.Range("K58").Value = 7576701
$ws.Range("L58").Value = 2303.6
$ws.Range("M58").Value = -7576498
$ws.Range("N58").Value = -2709.6
$ws.Range("H130").Value = 52520
$ws.Range("J130").Value = 52520
$ws.Range("L130").Value = 52520
$ws.Range("N130").Value = -62560
$ws.Range("H132").Value = 4763946
$ws.Range("I132").Value = 7144788.5
$ws.Range("J132").Value = 2259.9285
$ws.Range("K132").Value = 21434365.5
$ws.Range("L132").Value = 6779.7855
$ws.Range("M132").Value = -21431835.5
$ws.Range("N132").Value = -11839.7855
$ws.Range("H134").Value = 25646010
$ws.Range("I134").Value = 41673504
$ws.Range("J134").Value = 2022.8
$ws.Range("K134").Value = 125020512
$ws.Range("L134").Value = 6068.4
$ws.Range("M134").Value = -125017977
$ws.Range("N134").Value = -11138.4
$ws.Range("H136").Value = 5209702
$ws.Range("I136").Value = 7576701
$ws.Range("J136").Value = 2303.6
$ws.Range("K136").Value = 22730103
$ws.Range("L136").Value = 6910.799999999999
$ws.Range("M136").Value = -22727553
$ws.Range("N136").Value = -12010.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4661.2915
$ws.Range("J139").Value = 2896.8948
$ws.Range("L139").Value = 8690.6844
$ws.Range("N139").Value = -18970.6844

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5746.8057
$ws.Range("I70").Value = 5784.2964
$ws.Range("J70").Value = 5634.3335
$ws.Range("K70").Value = 5784.2964
$ws.Range("L70").Value = 5634.3335
$ws.Range("M70").Value = -5514.2964
$ws.Range("N70").Value = -6174.3335
$ws.Range("H73").Value = 5746.8057
$ws.Range("I73").Value = 5784.2964
$ws.Range("J73").Value = 5634.3335
$ws.Range("K73").Value = 5784.2964
$ws.Range("L73").Value = 5634.3335
$ws.Range("M73").Value = -4848.2964
$ws.Range("N73").Value = -7506.3335
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -5244

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 7000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 7000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 7000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -7876
$ws.Range("H45").Value = 14266.667
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 14266.667
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 14266.667
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -15080.667
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H132").Value = 13894319
$ws.Range("I132").Value = 16981120
$ws.Range("J132").Value = 3716.1667
$ws.Range("K132").Value = 50943360
$ws.Range("L132").Value = 11148.5001
$ws.Range("M132").Value = -50940830
$ws.Range("N132").Value = -16208.5001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69760
$ws.Range("J46").Value = 69760
$ws.Range("L46").Value = 69760
$ws.Range("N46").Value = -70222
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1300.2368
$ws.Range("I132").Value = 1003.3125
$ws.Range("K132").Value = 3009.9375
$ws.Range("M132").Value = -479.9375
$ws.Range("H134").Value = 69760
$ws.Range("J134").Value = 69760
$ws.Range("L134").Value = 209280
$ws.Range("N134").Value = -214350
$ws.Range("H136").Value = 3624377.8
$ws.Range("I136").Value = 696.9643
$ws.Range("J136").Value = 9261215
$ws.Range("K136").Value = 2090.8929
$ws.Range("L136").Value = 27783645
$ws.Range("M136").Value = 459.1071000000002
$ws.Range("N136").Value = -27788745
